$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Create Sheet2 right after Sheet1 and copy the original raw data
#    (the untouched Sheet1 content) into it.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row: reuse Sheet1's existing header cells/format (B1:E1 -> A1:D1)
$ws1.Range("B1:E1").Copy($ws2.Range("A1:D1"))
# Species header (originally F1) goes to E1 on Sheet2
$ws1.Range("F1").Copy($ws2.Range("E1"))

# Data rows 2-7: numeric columns B:E -> A:D (unformatted original data)
$ws1.Range("B2:E7").Copy($ws2.Range("A2:D7"))
# Species values (originally F2:F7) go to E2:E7
$ws1.Range("F2:F7").Copy($ws2.Range("E2:E7"))

# ---------------------------------------------------------------------
# 2) Remove the now-redundant Species column (F) from Sheet1 and
#    replace the row-index column (A) + numeric data (B:E) with the
#    summary statistics (pandas .describe()) of the four numeric
#    columns.
# ---------------------------------------------------------------------
$ws1.Columns.Item(6).Delete()

# Row labels. Plain words can go straight through .Value; the
# percent-looking labels ("25%", "50%", "75%") would otherwise be
# auto-converted to a numeric percentage by the normal cell-input
# parser, so those are produced as formula text in an unused scratch
# cell and brought across with a values-only paste, which preserves
# them as literal text without touching the style table.
$ws1.Cells.Item(2, 1).Value = "count"
$ws1.Cells.Item(3, 1).Value = "mean"
$ws1.Cells.Item(4, 1).Value = "std"
$ws1.Cells.Item(5, 1).Value = "min"

$percentLabels = @(
    @(6, "25%"),
    @(7, "50%"),
    @(8, "75%")
)
foreach ($pair in $percentLabels) {
    $r = $pair[0]
    $text = $pair[1]
    $ws1.Range("Z1").Formula = '="' + $text + '"'
    $ws1.Range("Z1").Copy()
    $ws1.Cells.Item($r, 1).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
    $ws1.Range("Z1").ClearContents()
}

$ws1.Cells.Item(9, 1).Value = "max"

$stats = @(
    @(6, 6, 6, 6),
    @(4.949999999999999, 3.383333333333333, 1.45, 0.2333333333333333),
    @(0.2880972058177588, 0.3430257521916782, 0.1378404875209022, 0.08164965809277261),
    @(4.6, 3, 1.3, 0.2),
    @(4.75, 3.125, 1.4, 0.2),
    @(4.95, 3.35, 1.4, 0.2),
    @(5.074999999999999, 3.575, 1.475, 0.2),
    @(5.4, 3.9, 1.7, 0.4)
)

for ($r = 0; $r -lt $stats.Length; $r++) {
    $row = $stats[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}

# Apply the bold/bordered header style (same as B1:E1) to the new
# label column A2:A9.
$ws1.Range("B1").Copy()
$ws1.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep Sheet1 as the active/selected sheet, matching the original file.
$ws1.Activate()
